# Applies the 15.7.1.1 workbook edit: adds a new "2023" column (N) to the
# table, mirroring the style/format of the existing neighbouring cells, and
# tweaks a couple of row heights / the data column width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row heights -----------------------------------------------------
# Row 1 grows a bit taller (wrapped header text), rows 2 & 3 get an
# explicit custom height.
$ws.Rows.Item(1).RowHeight = 55.5
$ws.Rows.Item(2).RowHeight = 13.5
$ws.Rows.Item(3).RowHeight = 13.5

# --- Column width ------------------------------------------------------
# Columns A:C (the label columns) get very slightly wider. The interop
# layer quantizes column widths to whole pixels, so 36.3 is the closest
# achievable setting to the target OOXML width of 37.140625 characters.
$ws.Range("A1:C1").ColumnWidth = 36.3

# --- New column N: header (row 4, "2023") ------------------------------
# Copy the format from the existing M4 (2022) header cell, then set value.
$ws.Range("M4").Copy() | Out-Null
$ws.Range("N4").PasteSpecial(-4122) | Out-Null
$ws.Range("N4").Value = 2023

# --- New column N: empty styled cell on row 3 ---------------------------
$ws.Range("M3").Copy() | Out-Null
$ws.Range("N3").PasteSpecial(-4122) | Out-Null

# --- New column N: data rows 5 & 6 (numeric 4.3499999999999996) --------
$ws.Range("J5").Copy() | Out-Null
$ws.Range("N5").PasteSpecial(-4122) | Out-Null
$ws.Range("N5").Value = 4.3499999999999996

$ws.Range("M6").Copy() | Out-Null
$ws.Range("N6").PasteSpecial(-4122) | Out-Null
$ws.Range("N6").Value = 4.3499999999999996

# --- New column N: row 7 gets "-" in a brand-new style ------------------
# (right+vertically centered, Times New Roman 9, no border - this does not
# match any existing style so a new cellXfs entry is created)
$ws.Range("N7").Value = "-"
$ws.Range("N7").Font.Name = "Times New Roman"
$ws.Range("N7").Font.Size = 9
$ws.Range("N7").Font.Bold = $false
$ws.Range("N7").HorizontalAlignment = -4152
$ws.Range("N7").VerticalAlignment = -4108

# --- New column N: row 8 gets "-" matching the M8 style -----------------
$ws.Range("M8").Copy() | Out-Null
$ws.Range("N8").PasteSpecial(-4122) | Out-Null
$ws.Range("N8").Value = "-"

$excel.CutCopyMode = 0
